# Add "2022-Q3" data:
#  - insert a new worksheet named "2022-Q3" right after "总计", pushing the
#    other quarter sheets (2022-Q2, 2021-Q3, 2020-Q4) one tab to the right
#  - populate the new sheet with the same layout/format as the other
#    per-quarter fund-holding sheets, filled with the 2022-Q3 data
#  - update the "总计" overview sheet with a new summary row for 2022-Q3

$wb = $excel.ActiveWorkbook

# xlPasteValues / xlPasteFormats constants used below to convert a helper
# formula into a plain static (text) value, and to copy cell formatting.
$xlPasteValues = -4163
$xlPasteFormats = -4122

# ---------------------------------------------------------------------
# 1. Insert the new "2022-Q3" worksheet right after "总计" (position 2)
# ---------------------------------------------------------------------
$wsTotal = $wb.Worksheets.Item(1)
$wsQ3 = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $wsTotal)
$wsQ3.Name = "2022-Q3"

# ---------------------------------------------------------------------
# 2. Copy formatting (styles/borders) of an existing per-quarter sheet's
#    header+data block into the new sheet, then overwrite the values.
#    After the insert, "2022-Q2" (old sheet2.xml content) now sits at
#    position 3 and still carries the correct header formatting.
#    (A straight PasteSpecial() on a freshly-added sheet only carries
#    values across, not styles, so copy formats explicitly instead.)
# ---------------------------------------------------------------------
$wsQ2 = $wb.Worksheets.Item(3)
$wsQ2.Range("A1:H2").Copy()
$wsQ3.Range("A1").PasteSpecial($xlPasteFormats)

$wsQ3.Cells.Item(1,2).Value = "基金代码"
$wsQ3.Cells.Item(1,3).Value = "基金名称"
$wsQ3.Cells.Item(1,4).Value = "基金规模"
$wsQ3.Cells.Item(1,5).Value = "股票总仓位"
$wsQ3.Cells.Item(1,6).Value = "仓位占比"
$wsQ3.Cells.Item(1,7).Value = "持有市值(亿元)"
$wsQ3.Cells.Item(1,8).Value = "仓位排名"

$wsQ3.Cells.Item(2,1).Value = 0
$wsQ3.Cells.Item(2,3).Value = "金元顺安价值增长混合"
$wsQ3.Cells.Item(2,8).Value = 8

# B2/D2/E2/F2/G2 look numeric but must stay plain text, matching the
# other per-quarter sheets (inline/shared string, not a number). Enter
# each one as a quoted-text formula, then flatten it in place to a
# static value (one cell at a time - a multi-cell PasteSpecial here
# leaves the last cell's formula behind) so no residual formula or
# extra number-format style is left on the cell.
$q3TextCells = @(
    @{ Col = 2; Text = "620004" },
    @{ Col = 4; Text = "0.34" },
    @{ Col = 5; Text = "74.99" },
    @{ Col = 6; Text = "1.92" },
    @{ Col = 7; Text = "0.0065" }
)
foreach ($cellInfo in $q3TextCells) {
    $cell = $wsQ3.Cells.Item(2, $cellInfo.Col)
    $cell.Formula = "=""" + $cellInfo.Text + """"
    $cell.Copy()
    $cell.PasteSpecial($xlPasteValues)
}

# ---------------------------------------------------------------------
# 3. Update the "总计" sheet: insert a summary row for 2022-Q3 right
#    after the header, shifting the existing quarter rows down by one
#    and renumbering the index column. Row 4 (2020-Q4) had no row below
#    it to copy formatting from, so give the new row 5 the same "A"
#    column style as the row above it before the shift.
# ---------------------------------------------------------------------
$wsTotal.Cells.Item(4,1).Copy()
$wsTotal.Cells.Item(5,1).PasteSpecial($xlPasteFormats)

# row 5 (2020-Q4, was row 4): 3, 2020-Q4, 3, 0.87
$wsTotal.Cells.Item(5,1).Value = 3
$wsTotal.Cells.Item(5,2).Value = "2020-Q4"
$wsTotal.Cells.Item(5,3).Value = 3
$wsTotal.Cells.Item(5,4).Value = 0.87

# row 4 (2021-Q3, was row 3): 2, 2021-Q3, 1, 0.02
$wsTotal.Cells.Item(4,1).Value = 2
$wsTotal.Cells.Item(4,2).Value = "2021-Q3"
$wsTotal.Cells.Item(4,3).Value = 1
$wsTotal.Cells.Item(4,4).Value = 0.02

# row 3 (2022-Q2, was row 2): 1, 2022-Q2, 1, 0.01
$wsTotal.Cells.Item(3,1).Value = 1
$wsTotal.Cells.Item(3,2).Value = "2022-Q2"
$wsTotal.Cells.Item(3,3).Value = 1
$wsTotal.Cells.Item(3,4).Value = 0.01

# row 2 (new): 0, 2022-Q3, 1, 0.01
$wsTotal.Cells.Item(2,1).Value = 0
$wsTotal.Cells.Item(2,2).Value = "2022-Q3"
$wsTotal.Cells.Item(2,3).Value = 1
$wsTotal.Cells.Item(2,4).Value = 0.01
